# tuck.c white box testing finished
# Applies the "Integrated White Box Test.xlsx" edit:
#   - shipment.c (sheet2): rename the cell under test to process_shipments(),
#     add two more test rows (findDiversion / printDiversion), widen a few
#     columns, wrap E2:F2.
#   - truck.c (sheet3): fill in the whichTruck() test case details and add
#     five more truck-helper test rows (cmpTruck, checkWeight, addWeight,
#     checkVolume, addVolume), widen columns, wrap text + taller row 2.
#   - mapping.c (sheet1) becomes the active/selected tab instead of truck.c.

$wb = $excel.ActiveWorkbook

$wsMapping  = $wb.Worksheets.Item("mapping.c")
$wsShipment = $wb.Worksheets.Item("shipment.c")
$wsTruck    = $wb.Worksheets.Item("truck.c")

# ---------------------------------------------------------------------
# shipment.c (sheet2)
# ---------------------------------------------------------------------

# Row 2 now documents process_shipments() rather than whichTruck().
$wsShipment.Range("A2").Value = "void process_shipments(struct Truck* trucksPtr)"
$wsShipment.Range("B2").Value = "11-162"

# E2/F2 are left blank but get the wrap-text style applied.
$wsShipment.Range("E2:F2").WrapText = $true

# Two new test-case rows.
$wsShipment.Range("A3").Value = "void findDiversion(const int truck, const struct Point P)"
$wsShipment.Range("B3").Value = "165-188"

$wsShipment.Range("A4").Value = "void printDiversion(const struct Route* diversion)"
$wsShipment.Range("B4").Value = "190-273"

# Column widths.
$wsShipment.Columns.Item(3).ColumnWidth = 44.5924479166667
$wsShipment.Columns.Item(4).ColumnWidth = 15.5924479166667
$wsShipment.Columns.Item(5).ColumnWidth = 14.5924479166667

# Page orientation.
$wsShipment.PageSetup.Orientation = 1

# View: selection moves to C30 (tabSelected is no longer set on this sheet).
$wsShipment.Range("C30").Select()

# ---------------------------------------------------------------------
# truck.c (sheet3)
# ---------------------------------------------------------------------

# Row 2 (whichTruck) gets its Function / Lines of Code / Case Description / Input / Expected / Actual filled in.
$wsTruck.Range("A2").Value = "int whichTruck(const struct Map* baseMap, struct Point dest, struct Truck* trucks_ptr)"
$wsTruck.Range("B2").Value = "7-105"
$wsTruck.Range("C2").Value = "Test if the function selects the correct truck based on the closest distance to the destination and the available space in the trucks. Test when the closest distance is on the blue, green and yellow routes, and when the closest distance is the same for two or three of them. Also, test when there is no available truck."
$wsTruck.Range("D2").Value = "No input required"
$wsTruck.Range("E2").Value = "If truck 1 is selected, return 1;`nIf truck 2 is selected, return 2;`nIf truck 3 is selected, return 3;`nIf no truck is available, return 0;"
$wsTruck.Range("F2").Value = "If truck 1 is selected, return 1;`nIf truck 2 is selected, return 2;`nIf truck 3 is selected, return 3;`nIf no truck is available, return 0;"
$wsTruck.Range("E2:F2").WrapText = $true
$wsTruck.Rows.Item(2).RowHeight = 180

# Row 3: cmpTruck
$wsTruck.Range("A3").Value = "int cmpTruck(struct Truck* trucksPtr, int truckIdx1, int truckIdx2)"
$wsTruck.Range("B3").Value = "107-118"
$wsTruck.Range("C3").Value = "Test if the function properly compares the weigh, volume, and load between two trucks and return the one that is less full."
$wsTruck.Range("D3").Value = "No input required"
$wsTruck.Range("E3").Value = "Returns the truck that is less full"
$wsTruck.Range("F3").Value = "Returns the truck that is less full"
$wsTruck.Range("G3").Formula = '=IF(E3=F3, "PASS")'

# Row 4: checkWeight
$wsTruck.Range("A4").Value = "int checkWeight(struct Truck* truckPtr, int truckIdx, int weight)"
$wsTruck.Range("B4").Value = "120-122"
$wsTruck.Range("C4").Value = "Test if the function correctly checks  whether adding a new package to a given truck would exceed its weight limit or not."
$wsTruck.Range("D4").Value = "No input required"
$wsTruck.Range("E4").Value = "It passes if a truck has more room and fails if it doesn" + [char]0x2019 + "t"
$wsTruck.Range("F4").Value = "It passes if a truck has more room and fails if it doesn" + [char]0x2019 + "t"

# Row 5: addWeight
$wsTruck.Range("A5").Value = "void addWeight(struct Truck* truckPtr, int truckIdx, int weight)"
$wsTruck.Range("B5").Value = "124-126"
$wsTruck.Range("C5").Value = "Test if the function adds the weight correctly to the truck."
$wsTruck.Range("D5").Value = "No input required"
$wsTruck.Range("E5").Value = "Correctly adds the weight to the truck"
$wsTruck.Range("F5").Value = "Correctly adds the weight to the truck"

# Row 6: checkVolume
$wsTruck.Range("A6").Value = "int checkVolume(struct Truck* trucksPtr, int truckIdx, double size)"
$wsTruck.Range("B6").Value = "128-130"
$wsTruck.Range("C6").Value = "Test if the function correctly checks if the package would exceed truck volume limit."
$wsTruck.Range("D6").Value = "No input required"
$wsTruck.Range("E6").Value = "Return 1 if the new weight exceeds limit, and 0 if it doesn" + [char]0x2019 + "t"
$wsTruck.Range("F6").Value = "Return 1 if the new weight exceeds limit, and 0 if it doesn" + [char]0x2019 + "t"

# Row 7: addVolume
$wsTruck.Range("A7").Value = "void addVolume(struct Truck* trucksPtr, int truckIdx, double size)"
$wsTruck.Range("B7").Value = "132-134"
$wsTruck.Range("C7").Value = "Test if the function correctly adds volume to the truck"
$wsTruck.Range("D7").Value = "No input required"
$wsTruck.Range("E7").Value = "Truck's volume correctly updated"
$wsTruck.Range("F7").Value = "Truck's volume correctly updated"

# Fill the PASS/FAIL formula down for the newly added rows (creates the
# shared-formula group the same way dragging the fill handle would).
$wsTruck.Range("G4:G7").Formula = '=IF(E4=F4, "PASS")'
$wsTruck.Range("G6").Formula = '=IF(E6=F6, "PASS")'
$wsTruck.Range("G7").Formula = '=IF(E7=F7, "PASS")'

# Column widths.
$wsTruck.Columns.Item(1).ColumnWidth = 73.5924479166667
$wsTruck.Columns.Item(2).ColumnWidth = 11.5924479166667
$wsTruck.Columns.Item(3).ColumnWidth = 102.5924479166667
$wsTruck.Columns.Item(4).ColumnWidth = 14.8776041666667
$wsTruck.Columns.Item(5).ColumnWidth = 13.5924479166667

# View: truck.c is no longer the selected tab, selection moves to A13.
$wsTruck.Range("A13").Select()

# ---------------------------------------------------------------------
# mapping.c (sheet1) becomes the active/selected tab
# ---------------------------------------------------------------------

$wsMapping.Activate()
$wsMapping.Range("C1").Select()
$wsMapping.Application.ActiveWindow.ScrollRow = 1
$wsMapping.Range("C25").Select()
